$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text storage (several
# prices parse as valid numbers, e.g. "126.30" -> 126.3, which would drop
# the exact formatted text the source data uses). Force text format first.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.602.36'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '1.802.85'
$ws.Range("E3").Value = '  -0.72%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = '316.72'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '0.5359'
$ws.Range("E7").Value = '  -6.09%  '
$ws.Range("D8").Value = '0.3763'
$ws.Range("E8").Value = '  -2.62%  '
$ws.Range("D9").Value = '0.07503'
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = '42.37'
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D13").Value = '20.68'
$ws.Range("E13").Value = '  -2.68%  '
$ws.Range("D14").Value = '6.147'
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").Value = '7.382'
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("D16").Value = '1.800.36'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("E17").Value = '  -2.15%  '
$ws.Range("D18").Value = '0.00001065'
$ws.Range("E18").Value = '  -1.16%  '
$ws.Range("D19").Value = '0.06445'
$ws.Range("E19").Value = '  -0.44%  '
$ws.Range("E20").Value = '  +0.16%  '
$ws.Range("E21").Value = '  -0.49%  '
$ws.Range("D22").Value = '5.925'
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").Value = '28.625.58'
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("D25").Value = '2.099'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").Value = '158.32'
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("E27").Value = '  -2.31%  '
$ws.Range("D28").Value = '2.007.37'
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '2.353'
$ws.Range("E29").Value = '  -4.00%  '
$ws.Range("D30").Value = '122.88'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("D31").Value = '1.105'
$ws.Range("E31").Value = '  -5.40%  '
$ws.Range("D32").Value = '0.1052'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").Value = '5.644'
$ws.Range("E33").Value = '  -2.53%  '
$ws.Range("D34").Value = '3.683'
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '8.734'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").Value = '5.042'
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  -4.00%  '
$ws.Range("D41").Value = '0.6223'
$ws.Range("E41").Value = '  -3.13%  '
$ws.Range("D42").Value = '1.195'
$ws.Range("E42").Value = '  +3.45%  '
$ws.Range("D43").Value = '1.423'
$ws.Range("E43").Value = '  +3.49%  '
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = '13.25'
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D48").Value = '126.30'
$ws.Range("E48").Value = '  +3.27%  '
$ws.Range("D49").Value = '1.940'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").Value = '1.154'
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = '0.06883'
$ws.Range("E51").Value = '  +0.45%  '

# Rows 35/36 swap (Algorand <-> Hedera) with refreshed price/volume data.
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.06507'
$ws.Range("E35").Value = '  +6.70%  '
$ws.Range("B36").Value = 'Algorand'
$ws.Range("C36").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D36").Value = '0.2247'
$ws.Range("E36").Value = '  +3.81%  '

# Rows 46/47 swap (PancakeSwap <-> Decentraland) with refreshed price/volume data.
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.5849'
$ws.Range("E46").Value = '  -2.56%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '3.686'
$ws.Range("E47").Value = '  -0.59%  '
